$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (values are stored as text/percent-as-text, matching the
# source data feed format used throughout this sheet).
$updates = [ordered]@{
    "D2" = "287.05"
    "E2" = "2.65%"
    "D3" = "28.64"
    "E3" = "4.16%"
    "D4" = "5.051"
    "E4" = "4.51%"
    "D5" = "0.06658"
    "E5" = "4.25%"
    "E6" = "4.41%"
    "D7" = "3.395"
    "E7" = "2.03%"
    "D8" = "1.370"
    "E8" = "4.56%"
    "D9" = "0.9418"
    "E9" = "5.22%"
    "D10" = "0.1559"
    "E10" = "1.02%"
    "D11" = "0.06642"
    "E11" = "-1.30%"
    "D12" = "0.07641"
    "E12" = "2.02%"
    "D13" = "0.02955"
    "E13" = "0.23%"
    "D14" = "0.08999"
    "E14" = "0.00%"
    "D15" = "0.001591"
    "E15" = "1.42%"
    "D16" = "0.04492"
    "E16" = "2.10%"
    "D17" = "0.0006480"
    "E17" = "-0.94%"
    "D18" = "0.006343"
    "E18" = "5.71%"
    "E19" = "-1.16%"
    "D20" = "2.259"
    "E20" = "1.38%"
    "E21" = "2.25%"
    "D22" = "0.1298"
    "E22" = "-3.94%"
    "D23" = "4.085"
    "E23" = "4.64%"
    "D24" = "0.1552"
    "E24" = "3.21%"
    "E25" = "0.37%"
    "D26" = "0.004493"
    "E26" = "5.08%"
    "D27" = "0.0001250"
    "E27" = "5.97%"
    "D28" = "0.0001618"
    "E28" = "-2.16%"
    "D40" = "0.04206"
    "E40" = "3.43%"
    "D41" = "0.006768"
    "E41" = "2.30%"
    "D42" = "0.1257"
    "E42" = "-10.33%"
    "D43" = "0.002020"
    "E43" = "-2.38%"
    "D44" = "0.01233"
    "E44" = "12.11%"
    "D45" = "0.00005678"
    "E45" = "2.31%"
    "D47" = "0.01307"
    "E47" = "-29.34%"
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    # Force a Text number format so Excel keeps the written value as a literal
    # string (matching the original inline-string cell contents) instead of
    # silently re-interpreting it as a number/percentage.
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
}
